$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Colombia" and "Ecuador" swap places in the shared string table, so the two
# rows that reference those entries now render the other country's label:
# row 29 (shared string index 33) now reads "Colombia"; row 30 (index 34) now
# reads "Ecuador". Set the text directly so the engine records this relabel.
$ws.Range("A29").Value = "Colombia"
$ws.Range("A30").Value = "Ecuador"

# Refresh the daily COVID-19 figures (Casos totales, Nuevos casos, Casos
# activos, Recuperados, Casos criticos, Muertes hoy, Muertes) for each country
# row that changed in this update.
# Row 4 - Estados Unidos
$ws.Range("B4").Value = 2088493
$ws.Range("C4").Value = 22092
$ws.Range("D4").Value = 813663
$ws.Range("E4").Value = 1158812
$ws.Range("G4").Value = 888
$ws.Range("H4").Value = 116018
# Row 5 - Brasil
$ws.Range("B5").Value = 802828
$ws.Range("C5").Value = 27644
$ws.Range("E5").Value = 365216
$ws.Range("G5").Value = 1123
$ws.Range("H5").Value = 40920
# Row 11 - Peru
$ws.Range("B11").Value = 214788
$ws.Range("C11").Value = 5965
$ws.Range("D11").Value = 102429
$ws.Range("E11").Value = 106250
$ws.Range("G11").Value = 206
$ws.Range("H11").Value = 6109
# Row 12 - Alemania
$ws.Range("B12").Value = 186795
$ws.Range("C12").Value = 285
$ws.Range("E12").Value = 6744
$ws.Range("G12").Value = 7
$ws.Range("H12").Value = 8851
# Row 20 - Canada
$ws.Range("B20").Value = 97530
$ws.Range("C20").Value = 405
$ws.Range("D20").Value = 57658
$ws.Range("E20").Value = 31876
# Row 29 - Colombia
$ws.Range("B29").Value = 45212
$ws.Range("C29").Value = 1530
$ws.Range("D29").Value = 17790
$ws.Range("E29").Value = 25934
$ws.Range("G29").Value = 55
$ws.Range("H29").Value = 1488
# Row 30 - Ecuador
$ws.Range("B30").Value = 44440
$ws.Range("D30").Value = 21862
$ws.Range("E30").Value = 18858
$ws.Range("H30").Value = 3720
# Row 49 - Japon
$ws.Range("B49").Value = 17292
$ws.Range("C49").Value = 41
$ws.Range("D49").Value = 15383
$ws.Range("E49").Value = 989
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 920
# Row 55 - Nigeria
$ws.Range("B55").Value = 14554
$ws.Range("C55").Value = 681
$ws.Range("D55").Value = 4494
$ws.Range("E55").Value = 9673
$ws.Range("G55").Value = 5
$ws.Range("H55").Value = 387
# Row 66 - Noruega
$ws.Range("B66").Value = 8608
$ws.Range("C66").Value = 14
$ws.Range("E66").Value = 228
# Row 81 - Guinea
$ws.Range("B81").Value = 4372
$ws.Range("C81").Value = 114
$ws.Range("D81").Value = 3033
$ws.Range("E81").Value = 1316
# Row 91 - Bulgaria
$ws.Range("B91").Value = 3086
$ws.Range("C91").Value = 93
$ws.Range("D91").Value = 1688
$ws.Range("E91").Value = 1230
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 168
# Row 132 - Uruguay
$ws.Range("D132").Value = 772
$ws.Range("E132").Value = 52
# Row 138 - Cabo Verde
$ws.Range("E138").Value = 357
$ws.Range("G138").Value = 1
$ws.Range("H138").Value = 6
# Row 169 - Guyana
$ws.Range("B169").Value = 158
$ws.Range("C169").Value = 2
$ws.Range("E169").Value = 54
